$d = $word.ActiveDocument

# Locate the "Full-Stack Development and Data Engineering" paragraph under the
# Siege Analytics / PARTNER entry — the new bullet points are inserted right
# after it and before the existing "• Lead comprehensive research..." bullet.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Full-Stack Development and Data Engineering`r") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    Write-Output "ERROR: anchor paragraph 'Full-Stack Development and Data Engineering' not found"
} else {
    $lines = @(
        "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
        "• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times",
        "• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis",
        "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
    )

    $cur = $anchor
    foreach ($line in $lines) {
        $cur.Range.InsertParagraphAfter()
        $cur = $cur.Next()
        $cur.Range.Text = $line
    }

    Write-Output "Inserted $($lines.Count) new bullet paragraphs after 'Full-Stack Development and Data Engineering'"
}
